# Fix circular-reference-prone SUM ranges in the LOE sheet's summary rows
# (Management / Closeout rows 31-33) so they point at the actual data rows
# (E3:E30 / E3:E31 / E3:E32 / G3:G32) instead of ranges that looped back
# over the formula cells themselves (E17:E44 / E17:E46 / G17:G46).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOE")

$ws.Range("D31").Formula2 = "=ROUND(SUM(E3:E30)*`$C`$30,0)"
$ws.Range("E31").Formula2 = "=ROUND(SUM(E3:E30)*0.25,0)"

$ws.Range("D32").Formula2 = "=ROUND(SUM(E3:E31)*0.20,0)"
$ws.Range("E32").Formula2 = "=ROUND(SUM(E3:E31)*0.20,0)"

$ws.Range("G33").Formula2 = '=TEXT(SUM(G3:G32),"$#,##0")'

# E33's new formula is a bare SUM(...) (no ROUND/TEXT wrapper), and this
# engine (like Excel itself) auto-inherits the number format of the first
# referenced cell for such "plain reference" formulas, which would
# clobber the existing "TOTAL HOURS" totals-row styling (s=54). Re-apply
# the original formatting (copied from the untouched sibling cell D33,
# which carries the same style) after setting the formula so only the
# formula text changes, matching the source edit.
$ws.Range("E33").Formula2 = "=SUM(E3:E32)"
$ws.Range("D33").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wb.Save()
